$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 05:10"

# Row 26 (Pakistan)
$ws.Range("B26").Value = 321218
$ws.Range("C26").Value = 755
$ws.Range("D26").Value = 305395
$ws.Range("E26").Value = 9209
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 6614

# Row 31 (Belgica)
$ws.Range("B31").Value = 181511
$ws.Range("C31").Value = 8271
$ws.Range("D31").Value = 20587
$ws.Range("E31").Value = 150646
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 34
$ws.Range("H31").Value = 10278

# Row 44 (Kazajistan)
$ws.Range("B44").Value = 109094
$ws.Range("C44").Value = 110
$ws.Range("D44").Value = 104525
$ws.Range("E44").Value = 2801

# Row 153 (Belice)
$ws.Range("B153").Value = 2619
$ws.Range("C153").Value = 34
$ws.Range("D153").Value = 1596
$ws.Range("E153").Value = 983
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 40
